# Fix: "bieu do gantt chi mot doan" - add missing entries
$wb = $excel.ActiveWorkbook

$wsName = $wb.Worksheets.Item("Name")
$wsProject = $wb.Worksheets.Item("Project")

# On the "Project" sheet, add the two new process steps next to "Cắm LED" (row 6)
$wsProject.Range("C6").Value = "Hàn gá"
$wsProject.Range("D6").Value = "Hàn full"

# Leave the selection on the Project sheet at C7 (no longer the active tab)
$wsProject.Range("C7").Select()

# On the "Name" sheet, append the new worker "Hạnh" as a new row
$wsName.Range("A10").Value = "Hạnh"
$wsName.Range("A10").Select()

# Make the "Name" sheet the active tab
$wsName.Activate()
